$d = $word.ActiveDocument

# --- 1. Append "(This is a change – Version for branch alternate)" in dark
#        red after the first paragraph's existing text, with two extra
#        trailing spaces added to the original (unstyled) run first. ---

$p1 = $d.Paragraphs(1)
$r = $p1.Range

# Two trailing spaces stay part of the original (black) run.
$r.InsertAfter("  ")

# Each of the three segments below becomes its own run (matching how the
# original author's edit was incrementally typed) all in dark red
# (C00000).
$segStart = $r.End - 1
$r.InsertAfter("(This is a change – Ve")
$segEnd = $r.End - 1
$d.Range($segStart, $segEnd).Font.Color = 192

$segStart = $r.End - 1
$r.InsertAfter("rsion for branch alternate")
$segEnd = $r.End - 1
$d.Range($segStart, $segEnd).Font.Color = 192

$segStart = $r.End - 1
$r.InsertAfter(")")
$segEnd = $r.End - 1
$d.Range($segStart, $segEnd).Font.Color = 192

# --- 2. Add a new, empty paragraph after the very last paragraph of the
#        document, shaded light grey (F9F9F9). ---

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Style = $d.Styles("Normal")
$newPara.Shading.Texture = 0
$newPara.Shading.ForegroundPatternColor = -16777216
$newPara.Shading.BackgroundPatternColor = 16382457
